$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-16, columns A-D (B left blank where value is $null)
$data = @(
  @(1180128, 4,    2, 1),
  @(1190255, $null,3, 1),
  @(1180274, $null,1, 0),
  @(1180056, $null,3, 1),
  @(1180041, $null,5, 0),
  @(1180606, $null,3, 0),
  @(1180456, 0,    2, 0),
  @(2200022, 4,    0, 0),
  @(1180552, 2,    2, 0),
  @(1180207, 3,    1, 0),
  @(1180045, 4,    4, 0),
  @(1180212, 5,    2, 0),
  @(1180155, 7,    2, 0),
  @(1470343, 7,    3, 0),
  @(1190172, $null,4, 1)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    if ($null -eq $r[1]) {
        $ws.Cells.Item($row, 2).ClearContents()
    } else {
        $ws.Cells.Item($row, 2).Value = $r[1]
    }
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Remove row 17 entirely (dimension shrinks from D17 to D16)
$ws.Rows.Item(17).Delete()
